# Append the new "Eval(polynomial x)" section to the end of the document,
# matching the blank paragraph + bold heading + three body paragraphs that
# follow the existing "Degree(polynomial)" section.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xmlFrag = @"
<w:p $wNs><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p $wNs><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Eval(polynomial x)</w:t></w:r></w:p><w:p $wNs><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Polynomial is the polynomial you want to evaluate, k is what you want to evaluate it at can take in any list but will convert to sparse before evaluation</w:t></w:r></w:p><w:p $wNs><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">First checks if the polynomial is empty, if it is return 0, if it is not empty it will check if its dense if it is dense returns </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>(eval (to-sparse polynomial 0) k)</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, which converts it to sparse and runs the program again. </w:t></w:r></w:p><w:p $wNs><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">If the polynomial is sparse then applys an equation to all items in the list, this takes the k value and applys it into our polynomial, then it will run through the equation starting with </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>powers, then multiplication and finally addition, after it returns the int.</w:t></w:r></w:p>
"@

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML($xmlFrag) | Out-Null
